# Fix naive component forecaster bug - updates the K/J/.../B "diagonal"
# band of forecast-error values in rows 24-52 of Sheet1, matching the
# corrected model output (Presentation state 11.02).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24
$ws.Cells.Item(24, 11).Value = -19.18321641127324
# Row 25
$ws.Cells.Item(25, 10).Value = -17.46096375569671
$ws.Cells.Item(25, 11).Value = -3.862991559831116
# Row 26
$ws.Cells.Item(26, 9).Value = -19.04253003720006
$ws.Cells.Item(26, 10).Value = -5.444557841334467
$ws.Cells.Item(26, 11).Value = -2.749039957006937
# Row 27
$ws.Cells.Item(27, 8).Value = -17.48993877135166
$ws.Cells.Item(27, 9).Value = -3.891966575486066
$ws.Cells.Item(27, 10).Value = -1.196448691158537
$ws.Cells.Item(27, 11).Value = 0.1128003469303707
# Row 28
$ws.Cells.Item(28, 7).Value = -17.14253003720006
$ws.Cells.Item(28, 8).Value = -3.544557841334466
$ws.Cells.Item(28, 9).Value = -0.8490399570069369
$ws.Cells.Item(28, 10).Value = 0.4602090810819703
$ws.Cells.Item(28, 11).Value = 1.98156889852946
# Row 29
$ws.Cells.Item(29, 6).Value = -19.48929825540327
$ws.Cells.Item(29, 7).Value = -5.891326059537676
$ws.Cells.Item(29, 8).Value = -3.195808175210146
$ws.Cells.Item(29, 9).Value = -1.886559137121239
$ws.Cells.Item(29, 10).Value = -0.3651993196737491
$ws.Cells.Item(29, 11).Value = -6.130940961804114
# Row 30
$ws.Cells.Item(30, 5).Value = -17.10057817327245
$ws.Cells.Item(30, 6).Value = -3.502605977406859
$ws.Cells.Item(30, 7).Value = -0.8070880930793289
$ws.Cells.Item(30, 8).Value = 0.5021609450095783
$ws.Cells.Item(30, 9).Value = 2.023520762457068
$ws.Cells.Item(30, 10).Value = -3.742220879673297
$ws.Cells.Item(30, 11).Value = 3.165469857548658
# Row 31
$ws.Cells.Item(31, 4).Value = -18.37606615945818
$ws.Cells.Item(31, 5).Value = -4.778093963592582
$ws.Cells.Item(31, 6).Value = -2.082576079265053
$ws.Cells.Item(31, 7).Value = -0.7733270411761453
$ws.Cells.Item(31, 8).Value = 0.7480327762713443
$ws.Cells.Item(31, 9).Value = -5.017708865859021
$ws.Cells.Item(31, 10).Value = 1.889981871362934
$ws.Cells.Item(31, 11).Value = -0.5203221319907101
# Row 32
$ws.Cells.Item(32, 3).Value = -21.24253003720006
$ws.Cells.Item(32, 4).Value = -7.644557841334466
$ws.Cells.Item(32, 5).Value = -4.949039957006937
$ws.Cells.Item(32, 6).Value = -3.63979091891803
$ws.Cells.Item(32, 7).Value = -2.11843110147054
$ws.Cells.Item(32, 8).Value = -7.884172743600905
$ws.Cells.Item(32, 9).Value = -0.9764820063789501
$ws.Cells.Item(32, 10).Value = -3.386786009732595
$ws.Cells.Item(32, 11).Value = -5.082042055103904
# Row 33
$ws.Cells.Item(33, 2).Value = -38.47883389241241
$ws.Cells.Item(33, 3).Value = -24.88086169654682
$ws.Cells.Item(33, 4).Value = -22.18534381221929
$ws.Cells.Item(33, 5).Value = -20.87609477413038
$ws.Cells.Item(33, 6).Value = -19.35473495668289
$ws.Cells.Item(33, 7).Value = -25.12047659881326
$ws.Cells.Item(33, 8).Value = -18.2127858615913
$ws.Cells.Item(33, 9).Value = -20.62308986494494
$ws.Cells.Item(33, 10).Value = -22.31834591031625
$ws.Cells.Item(33, 11).Value = -19.37902846813061
# Row 34
$ws.Cells.Item(34, 2).Value = 13.59797219586559
$ws.Cells.Item(34, 3).Value = 16.29349008019312
$ws.Cells.Item(34, 4).Value = 17.60273911828203
$ws.Cells.Item(34, 5).Value = 19.12409893572952
$ws.Cells.Item(34, 6).Value = 13.35835729359916
$ws.Cells.Item(34, 7).Value = 20.26604803082111
$ws.Cells.Item(34, 8).Value = 17.85574402746747
$ws.Cells.Item(34, 9).Value = 16.16048798209616
$ws.Cells.Item(34, 10).Value = 19.0998054242818
$ws.Cells.Item(34, 11).Value = 17.72090440315328
# Row 35
$ws.Cells.Item(35, 2).Value = 2.695517884327529
$ws.Cells.Item(35, 3).Value = 4.004766922416437
$ws.Cells.Item(35, 4).Value = 5.526126739863926
$ws.Cells.Item(35, 5).Value = -0.2396149022664389
$ws.Cells.Item(35, 6).Value = 6.668075834955516
$ws.Cells.Item(35, 7).Value = 4.257771831601872
$ws.Cells.Item(35, 8).Value = 2.562515786230562
$ws.Cells.Item(35, 9).Value = 5.501833228416203
$ws.Cells.Item(35, 10).Value = 4.12293220728769
$ws.Cells.Item(35, 11).Value = 5.627651391037844
# Row 36
$ws.Cells.Item(36, 2).Value = 1.309249038088907
$ws.Cells.Item(36, 3).Value = 2.830608855536397
$ws.Cells.Item(36, 4).Value = -2.935132786593968
$ws.Cells.Item(36, 5).Value = 3.972557950627987
$ws.Cells.Item(36, 6).Value = 1.562253947274342
$ws.Cells.Item(36, 7).Value = -0.1330020980969671
$ws.Cells.Item(36, 8).Value = 2.806315344088674
$ws.Cells.Item(36, 9).Value = 1.42741432296016
$ws.Cells.Item(36, 10).Value = 2.932133506710315
$ws.Cells.Item(36, 11).Value = 2.605412008371104
# Row 37
$ws.Cells.Item(37, 2).Value = 1.52135981744749
$ws.Cells.Item(37, 3).Value = -4.244381824682876
$ws.Cells.Item(37, 4).Value = 2.66330891253908
$ws.Cells.Item(37, 5).Value = 0.2530049091854352
$ws.Cells.Item(37, 6).Value = -1.442251136185874
$ws.Cells.Item(37, 7).Value = 1.497066305999766
$ws.Cells.Item(37, 8).Value = 0.1181652848712531
$ws.Cells.Item(37, 9).Value = 1.622884468621408
$ws.Cells.Item(37, 10).Value = 1.296162970282197
$ws.Cells.Item(37, 11).Value = 2.134262027593551
# Row 38
$ws.Cells.Item(38, 2).Value = -5.765741642130365
$ws.Cells.Item(38, 3).Value = 1.14194909509159
$ws.Cells.Item(38, 4).Value = -1.268354908262054
$ws.Cells.Item(38, 5).Value = -2.963610953633364
$ws.Cells.Item(38, 6).Value = -0.0242935114477234
$ws.Cells.Item(38, 7).Value = -1.403194532576236
$ws.Cells.Item(38, 8).Value = 0.1015246511739181
$ws.Cells.Item(38, 9).Value = -0.225196847165293
$ws.Cells.Item(38, 10).Value = 0.6129022101460611
$ws.Cells.Item(38, 11).Value = -2.06679770341519
# Row 39
$ws.Cells.Item(39, 2).Value = 6.907690737221955
$ws.Cells.Item(39, 3).Value = 4.497386733868311
$ws.Cells.Item(39, 4).Value = 2.802130688497001
$ws.Cells.Item(39, 5).Value = 5.741448130682642
$ws.Cells.Item(39, 6).Value = 4.362547109554129
$ws.Cells.Item(39, 7).Value = 5.867266293304283
$ws.Cells.Item(39, 8).Value = 5.540544794965072
$ws.Cells.Item(39, 9).Value = 6.378643852276426
$ws.Cells.Item(39, 10).Value = 3.698943938715175
$ws.Cells.Item(39, 11).Value = 4.934366632778122
# Row 40
$ws.Cells.Item(40, 2).Value = -2.410304003353644
$ws.Cells.Item(40, 3).Value = -4.105560048724954
$ws.Cells.Item(40, 4).Value = -1.166242606539313
$ws.Cells.Item(40, 5).Value = -2.545143627667827
$ws.Cells.Item(40, 6).Value = -1.040424443917672
$ws.Cells.Item(40, 7).Value = -1.367145942256883
$ws.Cells.Item(40, 8).Value = -0.529046884945529
$ws.Cells.Item(40, 9).Value = -3.20874679850678
$ws.Cells.Item(40, 10).Value = -1.973324104443833
$ws.Cells.Item(40, 11).Value = -0.2490198157013881
# Row 41
$ws.Cells.Item(41, 2).Value = -1.695256045371309
$ws.Cells.Item(41, 3).Value = 1.244061396814331
$ws.Cells.Item(41, 4).Value = -0.1348396243141821
$ws.Cells.Item(41, 5).Value = 1.369879559435973
$ws.Cells.Item(41, 6).Value = 1.043158061096761
$ws.Cells.Item(41, 7).Value = 1.881257118408115
$ws.Cells.Item(41, 8).Value = -0.7984427951531354
$ws.Cells.Item(41, 9).Value = 0.4369798989098115
$ws.Cells.Item(41, 10).Value = 2.161284187652257
$ws.Cells.Item(41, 11).Value = 2.527729380265284
# Row 42
$ws.Cells.Item(42, 2).Value = 2.939317442185641
$ws.Cells.Item(42, 3).Value = 1.560416421057127
$ws.Cells.Item(42, 4).Value = 3.065135604807282
$ws.Cells.Item(42, 5).Value = 2.738414106468071
$ws.Cells.Item(42, 6).Value = 3.576513163779425
$ws.Cells.Item(42, 7).Value = 0.896813250218174
$ws.Cells.Item(42, 8).Value = 2.132235944281121
$ws.Cells.Item(42, 9).Value = 3.856540233023566
$ws.Cells.Item(42, 10).Value = 4.222985425636594
$ws.Cells.Item(42, 11).Value = -1.252462014395462
# Row 43
$ws.Cells.Item(43, 2).Value = -1.378901021128513
$ws.Cells.Item(43, 3).Value = 0.1258181626216415
$ws.Cells.Item(43, 4).Value = -0.2009033357175696
$ws.Cells.Item(43, 5).Value = 0.6371957215937845
$ws.Cells.Item(43, 6).Value = -2.042504191967466
$ws.Cells.Item(43, 7).Value = -0.8070814979045196
$ws.Cells.Item(43, 8).Value = 0.9172227908379254
$ws.Cells.Item(43, 9).Value = 1.283667983450953
$ws.Cells.Item(43, 10).Value = -4.191779456581102
$ws.Cells.Item(43, 11).Value = -0.8310663749237079
# Row 44
$ws.Cells.Item(44, 2).Value = 1.504719183750155
$ws.Cells.Item(44, 3).Value = 1.177997685410944
$ws.Cells.Item(44, 4).Value = 2.016096742722298
$ws.Cells.Item(44, 5).Value = -0.6636031708389534
$ws.Cells.Item(44, 6).Value = 0.5718195232239935
$ws.Cells.Item(44, 7).Value = 2.296123811966439
$ws.Cells.Item(44, 8).Value = 2.662569004579467
$ws.Cells.Item(44, 9).Value = -2.81287843545259
$ws.Cells.Item(44, 10).Value = 0.5478346462048052
# Row 45
$ws.Cells.Item(45, 2).Value = -0.326721498339211
$ws.Cells.Item(45, 3).Value = 0.511377558972143
$ws.Cells.Item(45, 4).Value = -2.168322354589108
$ws.Cells.Item(45, 5).Value = -0.932899660526161
$ws.Cells.Item(45, 6).Value = 0.7914046282162839
$ws.Cells.Item(45, 7).Value = 1.157849820829312
$ws.Cells.Item(45, 8).Value = -4.317597619202744
$ws.Cells.Item(45, 9).Value = -0.9568845375453494
# Row 46
$ws.Cells.Item(46, 2).Value = 0.838099057311354
$ws.Cells.Item(46, 3).Value = -1.841600856249897
$ws.Cells.Item(46, 4).Value = -0.60617816218695
$ws.Cells.Item(46, 5).Value = 1.118126126555495
$ws.Cells.Item(46, 6).Value = 1.484571319168523
$ws.Cells.Item(46, 7).Value = -3.990876120863533
$ws.Cells.Item(46, 8).Value = -0.6301630392061384
# Row 47
$ws.Cells.Item(47, 2).Value = -2.679699913561251
$ws.Cells.Item(47, 3).Value = -1.444277219498304
$ws.Cells.Item(47, 4).Value = 0.2800270692441409
$ws.Cells.Item(47, 5).Value = 0.646472261857169
$ws.Cells.Item(47, 6).Value = -4.828975178174887
$ws.Cells.Item(47, 7).Value = -1.468262096517492
# Row 48
$ws.Cells.Item(48, 2).Value = 1.235422694062947
$ws.Cells.Item(48, 3).Value = 2.959726982805392
$ws.Cells.Item(48, 4).Value = 3.32617217541842
$ws.Cells.Item(48, 5).Value = -2.149275264613636
$ws.Cells.Item(48, 6).Value = 1.211437817043759
# Row 49
$ws.Cells.Item(49, 2).Value = 1.724304288742445
$ws.Cells.Item(49, 3).Value = 2.090749481355473
$ws.Cells.Item(49, 4).Value = -3.384697958676583
$ws.Cells.Item(49, 5).Value = -0.02398487701918839
# Row 50
$ws.Cells.Item(50, 2).Value = 0.3664451926130281
$ws.Cells.Item(50, 3).Value = -5.109002247419028
$ws.Cells.Item(50, 4).Value = -1.748289165761633
# Row 51
$ws.Cells.Item(51, 2).Value = -5.475447440032056
$ws.Cells.Item(51, 3).Value = -2.114734358374661
# Row 52
$ws.Cells.Item(52, 2).Value = 3.360713081657395
